# msz - first smoke test is running
#
# Adds a new "Truck Insurance" smoke-test row to the vehicle-insurance flow
# table on Tabelle1 (row 4), mirroring the existing "Automobile" smoke-test
# row, then refreshes the column widths and active selection the way Excel
# does after such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 4) under the existing Record/Process table.
$ws.Range("A4").Value = "103_TruckInsuranceAutomobile_001_SmokeTest"
$ws.Range("B4").Value = "var103_TruckInsuranceAutomobile_001_SmokeTest"
$ws.Range("C4").Value = "Open Truck Insurance"
$ws.Range("E4").Value = "103_TruckInsuranceAutomobile_001_SmokeTest"

# Re-fit the columns that now hold the longer process/variable names.
$ws.Columns("D:F").AutoFit()

# Leave the selection where the author ended up after entering the data.
$ws.Range("D21").Select() | Out-Null
